# Atualização de bases das ligas, do dia: 03-04-2024 às 22:09
# Updates rows 118-121 (match ids 116-119) on the "Finland Veikkausliiga" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Finland Veikkausliiga")

# --- Row 118 ---
$ws.Range("B118").Value = 7719795
$ws.Range("F118").Value = "FC Haka"
$ws.Range("G118").Value = "VPS Vaasa"
$ws.Range("K118").Value = 2.375
$ws.Range("L118").Value = 3.3
$ws.Range("M118").Value = 2.7
$ws.Range("N118").Value = 2.3
$ws.Range("O118").Value = 3.3
$ws.Range("P118").Value = 2.75
$ws.Range("Q118").Value = -0.25
$ws.Range("R118").Value = 2.05
$ws.Range("S118").Value = 1.8
$ws.Range("U118").Value = 2.025
$ws.Range("V118").Value = 1.825

# --- Row 119 ---
$ws.Range("B119").Value = 7719688
$ws.Range("F119").Value = "FC Ilves"
$ws.Range("G119").Value = "FC Lahti"
$ws.Range("K119").Value = 1.727
$ws.Range("L119").Value = 3.6
$ws.Range("M119").Value = 4.2
$ws.Range("N119").Value = 1.65
$ws.Range("O119").Value = 3.75
$ws.Range("P119").Value = 4.75
$ws.Range("Q119").Value = -0.75
$ws.Range("R119").Value = 1.875
$ws.Range("S119").Value = 1.975
$ws.Range("U119").Value = 1.825
$ws.Range("V119").Value = 2.025

# --- Row 120 ---
$ws.Range("B120").Value = 7719689
$ws.Range("F120").Value = "SJK"
$ws.Range("G120").Value = "AC Oulu"
$ws.Range("K120").Value = 1.833
$ws.Range("L120").Value = 3.5
$ws.Range("M120").Value = 3.75
$ws.Range("N120").Value = 1.8
$ws.Range("O120").Value = 3.5
$ws.Range("P120").Value = 4
$ws.Range("Q120").Value = -0.5
$ws.Range("R120").Value = 1.85
$ws.Range("S120").Value = 2

# --- Row 121 ---
$ws.Range("R121").Value = 1.975
$ws.Range("S121").Value = 1.875
